$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.850.03"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.870.62"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'0.7355"
$ws.Range("E5").Value = "  -5.41%  "
$ws.Range("D6").Value = "'241.82"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.3158"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  -4.38%  "
$ws.Range("D10").Value = "'0.07094"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").Value = "'0.08393"
$ws.Range("E11").Value = "  -9.57%  "
$ws.Range("D12").Value = "'0.7517"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").Value = "'5.405"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "1.859.48"
$ws.Range("E14").Value = "  -4.53%  "
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "29.844.37"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "'6.051"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "'13.56"
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("D19").Value = "'243.03"
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("D20").Value = "'0.000007820"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "'0.9990"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "2.117.09"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "'7.905"
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'0.1569"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "'9.317"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").Value = "'164.03"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "'18.58"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").Value = "'2.017"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").Value = "'1.477"
$ws.Range("E30").Value = "  +3.65%  "
$ws.Range("D31").Value = "'4.618"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "'1.530"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").Value = "'4.304"
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("D34").Value = "'0.05332"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").Value = "'1.233"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "'0.7518"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'0.9991"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'2.698"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "'0.01947"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D41").Value = "'0.4467"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").Value = "1.107.31"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "'6.082"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'72.17"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").Value = "'0.8604"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'102.91"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'7.717"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").Value = "'1.840"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").Value = "'3.053"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").Value = "2.014.95"
$ws.Range("E51").Value = "  -2.35%  "
